$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.068.64'
$ws.Range('E2').Value = '  +5.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.920.68'
$ws.Range('E3').Value = '  +2.72%  '
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.18'
$ws.Range('E5').Value = '  +4.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5227'
$ws.Range('E7').Value = '  +3.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4101'
$ws.Range('E8').Value = '  +5.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08521'
$ws.Range('E9').Value = '  +2.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.129'
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.84'
$ws.Range('E11').Value = '  +2.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.43'
$ws.Range('E12').Value = '  +9.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.440'
$ws.Range('E13').Value = '  +3.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.926.33'
$ws.Range('E14').Value = '  +3.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.427'
$ws.Range('E15').Value = '  +2.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '95.53'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001114'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06683'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.39'
$ws.Range('E20').Value = '  +3.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.0000'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.022'
$ws.Range('E22').Value = '  +1.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.071.42'
$ws.Range('E23').Value = '  +5.42%  '
$ws.Range('E24').Value = '  +2.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.206'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.146.13'
$ws.Range('E26').Value = '  +2.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.31'
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.11'
$ws.Range('E28').Value = '  +2.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.453'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '129.38'
$ws.Range('E30').Value = '  +1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.079'
$ws.Range('E31').Value = '  +3.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1056'
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.052'
$ws.Range('E33').Value = '  +5.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.634'
$ws.Range('E34').Value = '  +0.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02491'
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06626'
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2219'
$ws.Range('E37').Value = '  +2.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.235'
$ws.Range('E38').Value = '  +4.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.192'
$ws.Range('E39').Value = '  +3.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.903'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6548'
$ws.Range('E41').Value = '  +2.83%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.65'
$ws.Range('E42').Value = '  +5.02%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.245'
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6162'
$ws.Range('E44').Value = '  +2.59%  '
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('E46').Value = '  +2.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.083'
$ws.Range('E47').Value = '  +3.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.247'
$ws.Range('E48').Value = '  +2.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.57'
$ws.Range('E49').Value = '  +1.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.166'
$ws.Range('E50').Value = '  +7.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.94'
$ws.Range('E51').Value = '  +4.48%  '
